# "added refs to class" - remove the example/reference row (row 2, columns B:F)
# and the stray numeric reference value in F7, then tidy up the view state
# (zoom + selection) to match where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the example reference content that was filled in for row 2
# (columns B through F) - only the category in column A stays.
$ws.Range("B2:F2").Clear()

# Clear the stray numeric reference that had been placed in F7.
$ws.Range("F7").Clear()

# Row 2 was tall to fit the wrapped reference text; let it size back down
# now that the content is gone.
$ws.Rows("2:2").AutoFit()

# Leave the view the way the author left it: zoomed out a bit, with the
# selection resting on C10.
$ws.Range("C10").Select()
$excel.ActiveWindow.Zoom = 85
